# Updated cryptos list on Sun May  7 15:41:17 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell($sheet, $row, $col, $text) {
    # Force the cell to remain a text value (matches the source file, where
    # every Price/Volume cell is stored as an inline/shared string, even
    # when it "looks" numeric) and then strip the number-format / quote
    # -prefix style residue so the cell keeps the default style index.
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# r, D (price), E (volume 1h)
$updates = @(
    @{ Row = 2;  D = "29.125.71";     E = "  +1.63%  " },
    @{ Row = 3;  D = "1.931.68";      E = "  +2.32%  " },
    @{ Row = 4;  E = "  +0.18%  " },
    @{ Row = 5;  D = "326.77";        E = "  +1.45%  " },
    @{ Row = 6;  D = "1.005";         E = "  +0.32%  " },
    @{ Row = 7;  D = "0.4614";        E = "  +1.13%  " },
    @{ Row = 8;  D = "0.3835";        E = "  +1.32%  " },
    @{ Row = 9;  D = "0.07773";       E = "  +1.31%  " },
    @{ Row = 10; D = "0.9812";        E = "  +2.49%  " },
    @{ Row = 11; D = "22.57";         E = "  +3.26%  " },
    @{ Row = 12; D = "1.929.35";      E = "  +2.14%  " },
    @{ Row = 13; D = "6.993";         E = "  +0.83%  " },
    @{ Row = 14; D = "5.707";         E = "  +1.38%  " },
    @{ Row = 15; D = "0.07073";       E = "  +0.83%  " },
    @{ Row = 16; D = "1.007";         E = "  +0.37%  " },
    @{ Row = 17; D = "84.58";         E = "  +2.64%  " },
    @{ Row = 18; D = "0.000009560";   E = "  +1.11%  " },
    @{ Row = 19; D = "16.79";         E = "  +1.37%  " },
    @{ Row = 20; D = "1.005";         E = "  +0.35%  " },
    @{ Row = 21; D = "29.124.09";     E = "  +1.66%  " },
    @{ Row = 22; D = "5.354";         E = "  +0.68%  " },
    @{ Row = 23; D = "10.98";         E = "  +1.43%  " },
    @{ Row = 24; D = "2.082";         E = "  +0.72%  " },
    @{ Row = 25; D = "157.97";        E = "  +1.85%  " },
    @{ Row = 26; D = "19.14";         E = "  +1.14%  " },
    @{ Row = 27; D = "5.678";         E = "  +1.58%  " },
    @{ Row = 28; D = "118.21";        E = "  +1.48%  " },
    @{ Row = 29; D = "1.856";         E = "  +2.53%  " },
    @{ Row = 30; D = "0.09347";       E = "  +1.43%  " },
    @{ Row = 31; D = "0.8670";        E = "  +3.30%  " },
    @{ Row = 32; D = "5.138";         E = "  +2.05%  " },
    @{ Row = 33; D = "1.254";         E = "  +1.20%  " },
    @{ Row = 34; D = "3.017";         E = "  -1.18%  " },
    @{ Row = 35; E = "  +1.90%  " },
    @{ Row = 36; D = "1.160";         E = "  +1.54%  " },
    @{ Row = 37; D = "1.005";         E = "  +0.28%  " },
    @{ Row = 38; D = "0.02052";       E = "  +1.46%  " },
    @{ Row = 39; D = "3.063";         E = "  +13.95%  " },
    @{ Row = 40; D = "7.565";         E = "  +2.00%  " },
    @{ Row = 41; D = "0.5521";        E = "  +1.27%  " },
    @{ Row = 42; D = "0.1758";        E = "  +1.02%  " },
    @{ Row = 43; D = "9.395";         E = "  +2.83%  " },
    @{ Row = 44; D = "0.000002852";   E = "  -2.93%  " },
    @{ Row = 45; D = "2.216";         E = "  +6.51%  " },
    @{ Row = 46; E = "  +1.42%  " },
    @{ Row = 49; D = "1.781";         E = "  +1.09%  " },
    @{ Row = 50; D = "110.40";        E = "  -0.06%  " },
    @{ Row = 51; D = "1.006";         E = "  +0.46%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('D')) {
        Set-TextCell $ws $r 4 $u.D
    }
    if ($u.ContainsKey('E')) {
        Set-TextCell $ws $r 5 $u.E
    }
}

# Rows 47 and 48 swap places: EnergySwap <-> Cronos (name/link/price/volume)
$ws.Cells.Item(47, 2).Value = "Cronos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws 47 4 "0.06931"
Set-TextCell $ws 47 5 "  +2.62%  "

$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws 48 4 "11.22"
Set-TextCell $ws 48 5 "  +0.82%  "
